# Update Camp-Fpr2 LR-pairs sheet with new TPM-derived values.
# - Adds "ECs" as a new sending/target cluster alongside the existing "FAPs" cluster.
# - Existing rows (FAPs -> FAPs/Inflammatory-Mac/Resolving-Mac) get refreshed numeric values.
# - New rows are added for FAPs -> ECs and ECs -> ECs/FAPs/Inflammatory-Mac/Resolving-Mac.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Camp/Fpr2)
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Camp"
$ws.Cells.Item(2,3).Value = "Fpr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.023517
$ws.Cells.Item(2,8).Value = 0.070551
$ws.Cells.Item(2,9).Value = 0.3003230928370446
$ws.Cells.Item(2,10).Value = 0.3003230928370446
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.050657
$ws.Cells.Item(2,14).Value = 0.151971
$ws.Cells.Item(2,15).Value = 0.002123644810648064
$ws.Cells.Item(2,16).Value = 0.002123644810648064
$ws.Cells.Item(2,17).Value = 0.001191300669
$ws.Cells.Item(2,18).Value = 0.010721706021
$ws.Cells.Item(2,19).Value = 0.0006377795776211664
$ws.Cells.Item(2,20).Value = 0.0006377795776211664

# Row 3: ECs -> FAPs (Camp/Fpr2)
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Camp"
$ws.Cells.Item(3,3).Value = "Fpr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.023517
$ws.Cells.Item(3,8).Value = 0.070551
$ws.Cells.Item(3,9).Value = 0.3003230928370446
$ws.Cells.Item(3,10).Value = 0.3003230928370446
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.006255666666667
$ws.Cells.Item(3,14).Value = 3.018767
$ws.Cells.Item(3,15).Value = 0.0421842909114609
$ws.Cells.Item(3,16).Value = 0.0421842909114609
$ws.Cells.Item(3,17).Value = 0.023664114513
$ws.Cells.Item(3,18).Value = 0.212977030617
$ws.Cells.Item(3,19).Value = 0.01266891671566757
$ws.Cells.Item(3,20).Value = 0.01266891671566757

# Row 4: ECs -> Inflammatory-Mac (Camp/Fpr2)
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Camp"
$ws.Cells.Item(4,3).Value = "Fpr2"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.023517
$ws.Cells.Item(4,8).Value = 0.070551
$ws.Cells.Item(4,9).Value = 0.3003230928370446
$ws.Cells.Item(4,10).Value = 0.3003230928370446
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 18.813815
$ws.Cells.Item(4,14).Value = 56.441445
$ws.Cells.Item(4,15).Value = 0.7887135162611822
$ws.Cells.Item(4,16).Value = 0.7887135162611822
$ws.Cells.Item(4,17).Value = 0.442444487355
$ws.Cells.Item(4,18).Value = 3.982000386195
$ws.Cells.Item(4,19).Value = 0.2368688825659389
$ws.Cells.Item(4,20).Value = 0.2368688825659389

# Row 5: ECs -> Resolving-Mac (Camp/Fpr2)
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Camp"
$ws.Cells.Item(5,3).Value = "Fpr2"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.023517
$ws.Cells.Item(5,8).Value = 0.070551
$ws.Cells.Item(5,9).Value = 0.3003230928370446
$ws.Cells.Item(5,10).Value = 0.3003230928370446
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.983073
$ws.Cells.Item(5,14).Value = 11.949219
$ws.Cells.Item(5,15).Value = 0.1669785480167087
$ws.Cells.Item(5,16).Value = 0.1669785480167088
$ws.Cells.Item(5,17).Value = 0.09366992774099998
$ws.Cells.Item(5,18).Value = 0.843029349669
$ws.Cells.Item(5,19).Value = 0.05014751397781692
$ws.Cells.Item(5,20).Value = 0.05014751397781693

# Row 6: FAPs -> ECs (Camp/Fpr2)
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Camp"
$ws.Cells.Item(6,3).Value = "Fpr2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.05478866666666667
$ws.Cells.Item(6,8).Value = 0.164366
$ws.Cells.Item(6,9).Value = 0.6996769071629554
$ws.Cells.Item(6,10).Value = 0.6996769071629554
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.050657
$ws.Cells.Item(6,14).Value = 0.151971
$ws.Cells.Item(6,15).Value = 0.002123644810648064
$ws.Cells.Item(6,16).Value = 0.002123644810648064
$ws.Cells.Item(6,17).Value = 0.002775429487333334
$ws.Cells.Item(6,18).Value = 0.024978865386
$ws.Cells.Item(6,19).Value = 0.001485865233026897
$ws.Cells.Item(6,20).Value = 0.001485865233026897

# Row 7: FAPs -> FAPs (Camp/Fpr2)
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Camp"
$ws.Cells.Item(7,3).Value = "Fpr2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.05478866666666667
$ws.Cells.Item(7,8).Value = 0.164366
$ws.Cells.Item(7,9).Value = 0.6996769071629554
$ws.Cells.Item(7,10).Value = 0.6996769071629554
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.006255666666667
$ws.Cells.Item(7,14).Value = 3.018767
$ws.Cells.Item(7,15).Value = 0.0421842909114609
$ws.Cells.Item(7,16).Value = 0.0421842909114609
$ws.Cells.Item(7,17).Value = 0.05513140630244445
$ws.Cells.Item(7,18).Value = 0.496182656722
$ws.Cells.Item(7,19).Value = 0.02951537419579333
$ws.Cells.Item(7,20).Value = 0.02951537419579333

# Row 8: FAPs -> Inflammatory-Mac (Camp/Fpr2)
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Camp"
$ws.Cells.Item(8,3).Value = "Fpr2"
$ws.Cells.Item(8,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.05478866666666667
$ws.Cells.Item(8,8).Value = 0.164366
$ws.Cells.Item(8,9).Value = 0.6996769071629554
$ws.Cells.Item(8,10).Value = 0.6996769071629554
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 18.813815
$ws.Cells.Item(8,14).Value = 56.441445
$ws.Cells.Item(8,15).Value = 0.7887135162611822
$ws.Cells.Item(8,16).Value = 0.7887135162611822
$ws.Cells.Item(8,17).Value = 1.030783838763333
$ws.Cells.Item(8,18).Value = 9.277054548870002
$ws.Cells.Item(8,19).Value = 0.5518446336952433
$ws.Cells.Item(8,20).Value = 0.5518446336952433

# Row 9: FAPs -> Resolving-Mac (Camp/Fpr2)
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Camp"
$ws.Cells.Item(9,3).Value = "Fpr2"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.05478866666666667
$ws.Cells.Item(9,8).Value = 0.164366
$ws.Cells.Item(9,9).Value = 0.6996769071629554
$ws.Cells.Item(9,10).Value = 0.6996769071629554
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.983073
$ws.Cells.Item(9,14).Value = 11.949219
$ws.Cells.Item(9,15).Value = 0.1669785480167087
$ws.Cells.Item(9,16).Value = 0.1669785480167088
$ws.Cells.Item(9,17).Value = 0.218227258906
$ws.Cells.Item(9,18).Value = 1.964045330154
$ws.Cells.Item(9,19).Value = 0.1168310340388918
$ws.Cells.Item(9,20).Value = 0.1168310340388918
